$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared date string "19/07/2023" (used by rows 1-7 in column B) is being
# renamed to "01/08/2023". Re-write those cells with a leading apostrophe so
# the engine keeps them as quoted text (preserves the existing "quotePrefix"
# cell style) instead of re-detecting the format.
for ($r = 1; $r -le 7; $r++) {
    $ws.Range("B$r").Value = "'01/08/2023"
}

# Rows 8, 10 and 11 previously referenced the older date "17/06/2023"; they
# now move to the same "01/08/2023" value as the rest of the column (this is
# the "primary key" change called out in the commit message). Row 9 keeps the
# old date and is intentionally left untouched.
$ws.Range("B8").Value = "'01/08/2023"
$ws.Range("B10").Value = "'01/08/2023"
$ws.Range("B11").Value = "'01/08/2023"

# Update the active selection to match the author's final cursor position.
$ws.Range("B11").Select()
